$wb = $excel.ActiveWorkbook

# --- Update "addListItem" sheet: rename the login username from "Userten" to "Usereleven" ---
$wsAdd = $wb.Worksheets.Item("addListItem")
$wsAdd.Range("A2").Value = "Usereleven"

# --- Update "createUser" sheet: bump the test user id from 1032 to 1033 ---
$wsCreate = $wb.Worksheets.Item("createUser")
$wsCreate.Range("A2").Value = 1033

# --- Make "addListItem" the active/selected sheet (was "createUser") ---
$wsAdd.Activate()
$wsAdd.Range("A2").Select()
